$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 60

$ws.Cells.Item($newRow, 1).Value = "CompaNanny"
$ws.Cells.Item($newRow, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($newRow, 3).Value = "VGO"

# Assign the date as a literal-text formula first so the engine doesn't
# auto-coerce the ISO-looking string into a date serial number, then
# collapse the formula down to its static text result (matching how the
# other rows in this column store their dates as plain text).
$ws.Cells.Item($newRow, 4).Value = "=""2024-09-02"""
$ws.Cells.Item($newRow, 4).Copy()
$ws.Cells.Item($newRow, 4).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 1
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
